# ------------------------------------------------------------------
# Add a "Network config" (WiFi join count / first-connect flag) flash
# layout block to Sheet1, matching a new flash sector
# (USER_PARAMETER_START_SECTOR_ADDRESS0 / ...ADDRESS1), and rename the
# old "蓝牙重新连接次数" field into "重新连接次数" (row 13) while adding a
# new "首次连接标志" field (row 14).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing field description (row 13) ---------------------
$ws.Range("D13").Value = "重新连接次数"
$ws.Range("E13").Value = "0x0"

# --- Add the new "Network config" flash-sector column (F) for the
#     already-existing rows 2-13 --------------------------------------
$ws.Range("F1").Value = ""
$ws.Range("F2").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F3").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F4").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F5").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F6").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F7").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F8").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F9").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F10").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F11").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F12").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F13").Value = "USER_PARAMETER_START_SECTOR_ADDRESS1"

# --- Add the new row 14 (wifi_join_cnt / first-connect flag) ---------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "wifi_join_cnt"
$ws.Range("C14").Value = "u8"
$ws.Range("D14").Value = "首次连接标志"
$ws.Range("E14").Value = "0x1"
$ws.Range("F14").Value = "USER_PARAMETER_START_SECTOR_ADDRESS1"

# --- Match formatting (box border, left/vcenter align) on every cell
#     in the new rows / new column by cloning an already-formatted
#     cell's look-and-feel ---------------------------------------------
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:F12").PasteSpecial(-4122) | Out-Null

$ws.Range("E13").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:F13").Copy() | Out-Null
$ws.Range("A14:F14").PasteSpecial(-4122) | Out-Null

# Normalise the E5:E6 merged-cell borders (it becomes a plain full box
# border, same look as every other bordered cell on the sheet)
$ws.Range("D5").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null

# --- Column F width ----------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 41.3

# --- Selection / active cell -------------------------------------------
$ws.Range("D15").Select() | Out-Null

Write-Host "Edit complete"
